# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# zh-cn and de-de worksheets to reflect a freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 2
$wsZhCn.Range("E2").Value = "2016-03-13 15:05:47"
$wsZhCn.Range("H2").Value = "2016-03-13 15:06:06"

# de-de sheet, row 2
$wsDeDe.Range("E2").Value = "2016-03-13 15:05:50"
$wsDeDe.Range("H2").Value = "2016-03-13 15:06:19"
